$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 457, shifting existing rows 457-554 down to 458-555.
$ws.Rows.Item(457).Insert()

# Populate the newly inserted row 457 with the new data record.
$ws.Cells.Item(457, 1).Value = 8
$ws.Cells.Item(457, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(457, 3).Value = "Coquimbo"
$ws.Cells.Item(457, 4).Value = 44637
$ws.Cells.Item(457, 5).Value = 4
$ws.Cells.Item(457, 6).Value = 100112024
$ws.Cells.Item(457, 7).Value = "Choclo"
$ws.Cells.Item(457, 8).Value = "Dulce o Americano"
$ws.Cells.Item(457, 9).Value = "Primera"
$ws.Cells.Item(457, 10).Value = 25000
$ws.Cells.Item(457, 11).Value = 190
$ws.Cells.Item(457, 12).Value = 200
$ws.Cells.Item(457, 13).Value = 195
$ws.Cells.Item(457, 14).Value = "$/unidad"
$ws.Cells.Item(457, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(457, 16).Value = 195
$ws.Cells.Item(457, 17).Value = 1
$ws.Cells.Item(457, 18).Value = "Hortaliza"

# Ensure the date cell keeps the same date style/number format as the rest of column D.
$ws.Cells.Item(457, 4).NumberFormat = $ws.Cells.Item(458, 4).NumberFormat
